# "thay doi 2 thu muc" (change 2 folders):
#  1) ppt/comments          -> delete the lone reviewer comment on slide 4
#  2) ppt/slideLayouts (+slideMasters) -> refresh the cached "today" date
#     field (datetimeFigureOut) on the master and every layout from
#     12/14/2019 to 12/18/2019

$p = $ppt.ActivePresentation

# --- 1) Remove the reviewer comment -----------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $comments = $slide.Comments
    for ($ci = $comments.Count; $ci -ge 1; $ci--) {
        $comments.Item($ci).Delete()
    }
}

# --- 2) Update the cached date placeholder text ------------------------
$newDate = "12/18/2019"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                if ($shp.HasTextFrame) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
